$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary header figures -------------------------------------
$ws.Range("E11").Value = 231325        # VALOR MORA
$ws.Range("C13").Value = 7             # Cant. Trabajadores
$ws.Range("F13").Value = 6             # Cant. Periodos

# --- Trim the worker detail table down from 10 rows to 7 ---------------
# Deleting 3 rows out of the existing data block (rows 16-18) shifts the
# remaining rows up and preserves the bottom-border formatting that lives
# on the last row of the table (it slides from row 25 to row 22).
$ws.Range("B16:B18").EntireRow.Delete()

# --- Write the new worker detail data (rows 16-22) ----------------------
$data = @(
    @("CC", "9100670",    "CARLOS AUGUSTO CARDONA RESTREPO",   "1910", 42000, 952968),
    @("CC", "73350017",   "RICARDO POLO PATERNINA",             "1911", 33125, 877803),
    @("CC", "73571489",   "LENIN JACOB IBAÑEZ PEREZ",           "2103", 3511,  877803),
    @("CC", "1052075224", "WILMER ALBERTO VARGAS ZAPATA",       "2103", 68000, 1700000),
    @("CC", "73577260",   "MARIO RAFAEL GREY RODRIGUEZ",        "2109", 1211,  908526),
    @("CC", "73122996",   "ALEJANDRO ISMAEL GUETTE SAAVEDRA",   "2110", 36341, 908526),
    @("CC", "9145170",    "HUGO ALFONSO MERCADO ZABALETA",      "2201", 47137, 1178421)
)

$row = 16
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $row = $row + 1
}
